# Auto-generated: apply Ixion_Profits value corrections per commit diff.
# Workbook sheets ALC..WVR correspond to the original single-sheet row numbers
# referenced in the diff (H/I/J/K/L/M/N = price & profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4011587
$ws.Range("J17").Value = 4114434
$ws.Range("L17").Value = 12343302
$ws.Range("N17").Value = -12343638

$ws.Range("H63").Value = 34271
$ws.Range("J63").Value = 34271
$ws.Range("L63").Value = 34271
$ws.Range("N63").Value = -35519

$ws.Range("H66").Value = 34271
$ws.Range("J66").Value = 34271
$ws.Range("L66").Value = 102813
$ws.Range("N66").Value = -109053

$ws.Range("H82").Value = 1398.6666
$ws.Range("I82").Value = 1398.6666
$ws.Range("K82").Value = 4195.9998
$ws.Range("M82").Value = -3789.9998

$ws.Range("H85").Value = 1398.6666
$ws.Range("I85").Value = 1398.6666
$ws.Range("K85").Value = 4195.9998
$ws.Range("M85").Value = -2791.9998

$ws.Range("H112").Value = 1786687.1
$ws.Range("J112").Value = 2041871
$ws.Range("L112").Value = 6125613
$ws.Range("N112").Value = -6127829

$ws.Range("H129").Value = 1060.723
$ws.Range("J129").Value = 1067.4127
$ws.Range("L129").Value = 3202.2381
$ws.Range("N129").Value = -13202.2381

$ws.Range("H132").Value = 774.65753
$ws.Range("I132").Value = 629.43286
$ws.Range("K132").Value = 1888.29858
$ws.Range("M132").Value = 641.7014199999999

$ws.Range("H138").Value = 1707.5571
$ws.Range("J138").Value = 2443.3513
$ws.Range("L138").Value = 7330.053899999999
$ws.Range("N138").Value = -17610.0539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1353.0278
$ws.Range("I74").Value = 1230.1305
$ws.Range("K74").Value = 1230.1305
$ws.Range("M74").Value = -356.1305

$ws.Range("H77").Value = 1353.0278
$ws.Range("I77").Value = 1230.1305
$ws.Range("K77").Value = 6150.6525
$ws.Range("M77").Value = -1782.6525

$ws.Range("H102").Value = 2470944
$ws.Range("I102").Value = 2470944
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2470944
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2469322
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 1283613.1
$ws.Range("I122").Value = 1283613.1
$ws.Range("K122").Value = 3850839.3
$ws.Range("M122").Value = -3848389.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 4355
$ws.Range("I23").Value = 1500
$ws.Range("K23").Value = 1500
$ws.Range("M23").Value = -1217

$ws.Range("H80").Value = 343.03845
$ws.Range("J80").Value = 466.44446
$ws.Range("L80").Value = 466.44446
$ws.Range("N80").Value = -2462.44446

$ws.Range("H83").Value = 343.03845
$ws.Range("J83").Value = 466.44446
$ws.Range("L83").Value = 2332.2223
$ws.Range("N83").Value = -12316.2223

$ws.Range("H86").Value = 1999.875
$ws.Range("J86").Value = 1999.6666
$ws.Range("L86").Value = 1999.6666
$ws.Range("N86").Value = -4245.6666

$ws.Range("H89").Value = 1999.875
$ws.Range("J89").Value = 1999.6666
$ws.Range("L89").Value = 9998.333
$ws.Range("N89").Value = -21230.333

$ws.Range("H99").Value = 71430050
$ws.Range("I99").Value = 100001410
$ws.Range("K99").Value = 100001410
$ws.Range("M99").Value = -99999912

$ws.Range("H105").Value = 2282.7856
$ws.Range("I105").Value = 2048.9
$ws.Range("J105").Value = 2867.5
$ws.Range("K105").Value = 2048.9
$ws.Range("L105").Value = 2867.5
$ws.Range("M105").Value = -301.9000000000001
$ws.Range("N105").Value = -6361.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3856.6323
$ws.Range("I31").Value = 1451.4324
$ws.Range("J31").Value = 6727.355
$ws.Range("K31").Value = 1451.4324
$ws.Range("L31").Value = 6727.355
$ws.Range("M31").Value = -1156.4324
$ws.Range("N31").Value = -7317.355

$ws.Range("H34").Value = 3856.6323
$ws.Range("I34").Value = 1451.4324
$ws.Range("J34").Value = 6727.355
$ws.Range("K34").Value = 1451.4324
$ws.Range("L34").Value = 6727.355
$ws.Range("M34").Value = -1249.4324
$ws.Range("N34").Value = -7131.355

$ws.Range("H94").Value = 3828.4688
$ws.Range("I94").Value = 2938.5334
$ws.Range("J94").Value = 4613.706
$ws.Range("K94").Value = 2938.5334
$ws.Range("L94").Value = 4613.706
$ws.Range("M94").Value = -2487.5334
$ws.Range("N94").Value = -5515.706

$ws.Range("H105").Value = 1503.5454
$ws.Range("I105").Value = 1487.7222
$ws.Range("J105").Value = 1574.75
$ws.Range("K105").Value = 1487.7222
$ws.Range("L105").Value = 1574.75
$ws.Range("M105").Value = 259.2778000000001
$ws.Range("N105").Value = -5068.75

$ws.Range("H132").Value = 3038.2856
$ws.Range("I132").Value = 3104
$ws.Range("J132").Value = 2989
$ws.Range("K132").Value = 9312
$ws.Range("L132").Value = 8967
$ws.Range("M132").Value = -6782
$ws.Range("N132").Value = -14027

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 649
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 2471.9092
$ws.Range("I64").Value = 1870.3334
$ws.Range("J64").Value = 2697.5
$ws.Range("K64").Value = 5611.0002
$ws.Range("L64").Value = 8092.5
$ws.Range("M64").Value = -5341.0002
$ws.Range("N64").Value = -8632.5

$ws.Range("H66").Value = 649
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 2471.9092
$ws.Range("I67").Value = 1870.3334
$ws.Range("J67").Value = 2697.5
$ws.Range("K67").Value = 5611.0002
$ws.Range("L67").Value = 8092.5
$ws.Range("M67").Value = -4675.0002
$ws.Range("N67").Value = -9964.5

$ws.Range("H87").Value = 5042.6665
$ws.Range("I87").Value = 5042.6665
$ws.Range("K87").Value = 15127.9995
$ws.Range("M87").Value = -13879.9995

$ws.Range("H90").Value = 5042.6665
$ws.Range("I90").Value = 5042.6665
$ws.Range("K90").Value = 45383.9985
$ws.Range("M90").Value = -39143.9985

$ws.Range("H113").Value = 625582.06
$ws.Range("I113").Value = 597.75
$ws.Range("J113").Value = 2500535
$ws.Range("K113").Value = 1793.25
$ws.Range("L113").Value = 7501605
$ws.Range("M113").Value = 376.75
$ws.Range("N113").Value = -7505945

$ws.Range("H121").Value = 19275.777
$ws.Range("I121").Value = 592.25
$ws.Range("J121").Value = 34222.6
$ws.Range("K121").Value = 1776.75
$ws.Range("L121").Value = 102667.8
$ws.Range("M121").Value = -466.75
$ws.Range("N121").Value = -105287.8

$ws.Range("H131").Value = 935.9
$ws.Range("I131").Value = 365.8
$ws.Range("J131").Value = 965.9053
$ws.Range("K131").Value = 1097.4
$ws.Range("L131").Value = 2897.7159
$ws.Range("M131").Value = 3942.6
$ws.Range("N131").Value = -12977.7159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 132.18182
$ws.Range("I2").Value = 96.5
$ws.Range("J2").Value = 227.33333
$ws.Range("K2").Value = 96.5
$ws.Range("L2").Value = 227.33333
$ws.Range("M2").Value = 16.5
$ws.Range("N2").Value = -453.33333

$ws.Range("H107").Value = 1390.5
$ws.Range("I107").Value = 654.26666
$ws.Range("J107").Value = 3599.2
$ws.Range("K107").Value = 654.26666
$ws.Range("L107").Value = 3599.2
$ws.Range("M107").Value = 1265.73334
$ws.Range("N107").Value = -7439.2

$ws.Range("H126").Value = 5614.6206
$ws.Range("I126").Value = 8757.714
$ws.Range("J126").Value = 2681.0667
$ws.Range("K126").Value = 26273.142
$ws.Range("L126").Value = 8043.2001
$ws.Range("M126").Value = -23803.142
$ws.Range("N126").Value = -12983.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 36304.758
$ws.Range("I7").Value = 43339.082
$ws.Range("K7").Value = 43339.082
$ws.Range("M7").Value = -43227.082

$ws.Range("H40").Value = 500500000
$ws.Range("I40").Value = 500500000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 500500000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -500499864
$ws.Range("N40").ClearContents()

$ws.Range("H82").Value = 208919.8
$ws.Range("I82").Value = 2800
$ws.Range("K82").Value = 2800
$ws.Range("M82").Value = -2439

$ws.Range("H85").Value = 208919.8
$ws.Range("I85").Value = 2800
$ws.Range("K85").Value = 2800
$ws.Range("M85").Value = -1552

$ws.Range("H100").Value = 1236.125
$ws.Range("I100").Value = 964.8333
$ws.Range("J100").Value = 2050
$ws.Range("K100").Value = 964.8333
$ws.Range("L100").Value = 2050
$ws.Range("M100").Value = -423.8333
$ws.Range("N100").Value = -3132

$ws.Range("H126").Value = 36304.758
$ws.Range("I126").Value = 43339.082
$ws.Range("K126").Value = 130017.246
$ws.Range("M126").Value = -127547.246

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2045.9286
$ws.Range("I132").Value = 1645.7142
$ws.Range("K132").Value = 4937.142599999999
$ws.Range("M132").Value = -2407.142599999999
